# Weekly update: add this week's "Repollo" (Vega Central Mapocho de Santiago)
# price rows. Three new observations (date 2022-07-05, serial 44747) are
# inserted right after the existing row 516, pushing the former rows
# 517:535 down to 520:538 (dimension grows from A1:R535 to A1:R538).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 517:518:519 - everything that was 517:535 shifts to 520:538.
$ws.Rows("517:519").Insert()

# Row 517 - Crespo record / Primera
$ws.Range("A517").Value = 9
$ws.Range("B517").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C517").Value = "Metropolitana"
$ws.Range("D517").Value = 44747
$ws.Range("E517").Value = 13
$ws.Range("F517").Value = 100112006
$ws.Range("G517").Value = "Repollo"
$ws.Range("H517").Value = "Crespo record"
$ws.Range("I517").Value = "Primera"
$ws.Range("J517").Value = 3400
$ws.Range("K517").Value = 1500
$ws.Range("L517").Value = 1600
$ws.Range("M517").Value = 1550
$ws.Range("N517").Value = "`$/unidad"
$ws.Range("O517").Value = "Región Metropolitana"
$ws.Range("P517").Value = 1550
$ws.Range("Q517").Value = 1
$ws.Range("R517").Value = "Hortaliza"

# Row 518 - Crespo record / Segunda
$ws.Range("A518").Value = 9
$ws.Range("B518").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C518").Value = "Metropolitana"
$ws.Range("D518").Value = 44747
$ws.Range("E518").Value = 13
$ws.Range("F518").Value = 100112006
$ws.Range("G518").Value = "Repollo"
$ws.Range("H518").Value = "Crespo record"
$ws.Range("I518").Value = "Segunda"
$ws.Range("J518").Value = 1060
$ws.Range("K518").Value = 1200
$ws.Range("L518").Value = 1200
$ws.Range("M518").Value = 1200
$ws.Range("N518").Value = "`$/unidad"
$ws.Range("O518").Value = "Región Metropolitana"
$ws.Range("P518").Value = 1200
$ws.Range("Q518").Value = 1
$ws.Range("R518").Value = "Hortaliza"

# Row 519 - Morada(o) / Primera
$ws.Range("A519").Value = 9
$ws.Range("B519").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C519").Value = "Metropolitana"
$ws.Range("D519").Value = 44747
$ws.Range("E519").Value = 13
$ws.Range("F519").Value = 100112006
$ws.Range("G519").Value = "Repollo"
$ws.Range("H519").Value = "Morada(o)"
$ws.Range("I519").Value = "Primera"
$ws.Range("J519").Value = 970
$ws.Range("K519").Value = 1700
$ws.Range("L519").Value = 1700
$ws.Range("M519").Value = 1700
$ws.Range("N519").Value = "`$/unidad"
$ws.Range("O519").Value = "Región Metropolitana"
$ws.Range("P519").Value = 1700
$ws.Range("Q519").Value = 1
$ws.Range("R519").Value = "Hortaliza"
